# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 6 de Julio de 2020 a las 10:42"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 2983142
$ws.Range("C4").Value = 214
$ws.Range("D4").Value = 1289687
$ws.Range("E4").Value = 1560884
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 132571

# Row 7 - Rusia
$ws.Range("B7").Value = 687862
$ws.Range("C7").Value = 6611
$ws.Range("D7").Value = 454329
$ws.Range("E7").Value = 223237
$ws.Range("G7").Value = 135
$ws.Range("H7").Value = 10296

# Row 40 - Singapur
$ws.Range("B40").Value = 44983
$ws.Range("C40").Value = 183
$ws.Range("E40").Value = 4516

# Row 46 - Polonia
$ws.Range("D46").Value = 23966
$ws.Range("E46").Value = 10467

# Row 49 - Israel
$ws.Range("B49").Value = 30162
$ws.Range("C49").Value = 204
$ws.Range("D49").Value = 17974
$ws.Range("E49").Value = 11856
$ws.Range("G49").Value = 1
$ws.Range("H49").Value = 332

# Row 114 - Estonia
$ws.Range("B114").Value = 1994
$ws.Range("C114").Value = 1
$ws.Range("D114").Value = 1875

# Row 117 - Lituania
$ws.Range("B117").Value = 1841
$ws.Range("C117").Value = 5
$ws.Range("D117").Value = 1547
$ws.Range("E117").Value = 215

# Row 119 - Eslovaquia
$ws.Range("B119").Value = 1765
$ws.Range("C119").Value = 1
$ws.Range("E119").Value = 271

# Row 120 - Eslovenia
$ws.Range("B120").Value = 1716
$ws.Range("C120").Value = 16
$ws.Range("E120").Value = 221
